$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" (summary) sheet: push the existing 2022-Q3 totals row down to
#    row 3 and insert a new row 2 for 2022-Q4.
# ---------------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item(1)
$wsTotal.Range("A2:D2").Copy($wsTotal.Range("A3:D3"))
$wsTotal.Range("A3").Value = 1
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.01

# ---------------------------------------------------------------------------
# 2) The existing "2022-Q3" worksheet (rId2 / sheetId 2) keeps its position
#    right after "总计", but becomes the new "2022-Q4" sheet with fresh
#    holdings data. Its old contents are copied first onto a brand-new
#    worksheet which takes over the "2022-Q3" name.
# ---------------------------------------------------------------------------
$wsQ4 = $wb.Worksheets.Item(2)

$wsQ3 = $wb.Worksheets.Add($null, $wsQ4)
$wsQ4.Range("A1:H2").Copy($wsQ3.Range("A1:H2"))

$wsQ4.Name = "2022-Q4"
$wsQ3.Name = "2022-Q3"

# ---------------------------------------------------------------------------
# 3) Overwrite "2022-Q4" with the new quarter's fund holdings (two rows
#    instead of one). Re-stamp the header / index-column style (matches the
#    "总计" sheet's bold-centered style) before writing the new values.
# ---------------------------------------------------------------------------
$wsTotal.Range("B1").Copy($wsQ4.Range("B1:H1"))
$wsTotal.Range("A2").Copy($wsQ4.Range("A2:A3"))

$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

$wsQ4.Range("A2").Value = 0
$wsQ4.Range("B2").NumberFormat = "@"
$wsQ4.Range("B2").Value = "003241"
$wsQ4.Range("C2").Value = "创金合信量化发现灵活配置混合A"
$wsQ4.Range("D2:G2").NumberFormat = "@"
$wsQ4.Range("D2").Value = "0.30"
$wsQ4.Range("E2").Value = "91.60"
$wsQ4.Range("F2").Value = "1.52"
$wsQ4.Range("G2").Value = "0.0046"
$wsQ4.Range("H2").Value = 2

$wsQ4.Range("A3").Value = 1
$wsQ4.Range("B3").NumberFormat = "@"
$wsQ4.Range("B3").Value = "003242"
$wsQ4.Range("C3").Value = "创金合信量化发现灵活配置混合C"
$wsQ4.Range("D3:G3").NumberFormat = "@"
$wsQ4.Range("D3").Value = "0.29"
$wsQ4.Range("E3").Value = "91.60"
$wsQ4.Range("F3").Value = "1.52"
$wsQ4.Range("G3").Value = "0.0044"
$wsQ4.Range("H3").Value = 2

# Restore "总计" as the active sheet, matching the original view.
$wsTotal.Activate()
